$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the existing table ("Tabela1") to include the new 4th column first,
# so the new column header cell (D1) gets associated with the table.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D4"))

# Fill in the new column's header and data (this also populates the shared
# strings table in the same order as the target file: header, rule about
# numbers, rule about special characters, rule about spaces).
$ws.Range("D1").Value = "Template_rules_1"
$ws.Range("D3").Value = "Do not use number in the beginning of the column name"
$ws.Range("D4").Value = "Do not use special characters in the column name"
$ws.Range("D2").Value = 'Do not use space bar in column name, use "_" to separate'

# Make sure the table header picks up the real column name instead of the
# generic "Column4" placeholder.
$lo.ListColumns.Item(4).Name = "Template_rules_1"

# Give column D a sensible custom width similar to the other bestFit columns.
$ws.Columns.Item(4).ColumnWidth = 52.28515625
